$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells that receive a plain decimal-looking string must be forced to text format
# first, otherwise Excel auto-converts them to numbers and the literal formatting
# (trailing zeros, thousand-group dots, etc.) would be lost.
$textCells = @(
    "D2", "D3", "D5", "D6", "D8", "D9", "D10", "D11", "D12", "D14", "D15", "D16", "D17", "D18", "D19", "D20", "D21", "D22", "D23", "D24", "D25", "D27", "D29", "D30", "D31", "D32", "D33", "D34", "D35", "D36", "D37", "D38", "D39", "D40", "D41", "D42", "D43", "D45", "D46", "D47", "D48", "D49", "D50", "D51"
)
foreach ($cellRef in $textCells) {
    $ws.Range($cellRef).NumberFormat = "@"
}

$ws.Range("D2").Value = "54.309.10"
$ws.Range("E2").Value = "  -6.62%  "
$ws.Range("D3").Value = "2.876.81"
$ws.Range("E3").Value = "  -9.72%  "
$ws.Range("E4").Value = "  -0.07%  "
$ws.Range("D5").Value = "468.82"
$ws.Range("E5").Value = "  -11.89%  "
$ws.Range("D6").Value = "125.36"
$ws.Range("E6").Value = "  -6.82%  "
$ws.Range("E7").Value = "  -0.13%  "
$ws.Range("D8").Value = "2.872.35"
$ws.Range("E8").Value = "  -9.90%  "
$ws.Range("D9").Value = "0.402"
$ws.Range("E9").Value = "  -11.40%  "
$ws.Range("D10").Value = "6.60"
$ws.Range("E10").Value = "  -9.44%  "
$ws.Range("D11").Value = "0.0953"
$ws.Range("E11").Value = "  -15.01%  "
$ws.Range("D12").Value = "0.326"
$ws.Range("E12").Value = "  -17.27%  "
$ws.Range("E13").Value = "  -5.02%  "
$ws.Range("D14").Value = "3.361.79"
$ws.Range("E14").Value = "  -9.97%  "
$ws.Range("D15").Value = "23.10"
$ws.Range("E15").Value = "  -10.01%  "
$ws.Range("D16").Value = "54.252.38"
$ws.Range("E16").Value = "  -6.95%  "
$ws.Range("D17").Value = "2.874.82"
$ws.Range("E17").Value = "  -9.84%  "
$ws.Range("D18").Value = "0.0000132"
$ws.Range("E18").Value = "  -14.83%  "
$ws.Range("D19").Value = "5.31"
$ws.Range("E19").Value = "  -9.27%  "
$ws.Range("D20").Value = "11.35"
$ws.Range("E20").Value = "  -14.18%  "
$ws.Range("D21").Value = "7.03"
$ws.Range("E21").Value = "  -13.31%  "
$ws.Range("D22").Value = "295.97"
$ws.Range("E22").Value = "  -17.71%  "
$ws.Range("D23").Value = "1.00"
$ws.Range("E23").Value = "  +0.14%  "
$ws.Range("D24").Value = "0.437"
$ws.Range("E24").Value = "  -15.28%  "
$ws.Range("D25").Value = "58.28"
$ws.Range("E25").Value = "  -16.43%  "
$ws.Range("E26").Value = "  +0.08%  "
$ws.Range("D27").Value = "0.150"
$ws.Range("E27").Value = "  -10.55%  "
$ws.Range("E28").Value = "  -0.02%  "
$ws.Range("D29").Value = "0.0₃0801"
$ws.Range("E29").Value = "  -15.82%  "
$ws.Range("D30").Value = "6.06"
$ws.Range("E30").Value = "  -12.34%  "
$ws.Range("D31").Value = "1.11"
$ws.Range("E31").Value = "  -8.03%  "
$ws.Range("D32").Value = "6.11"
$ws.Range("E32").Value = "  -12.30%  "
$ws.Range("D33").Value = "1.60"
$ws.Range("E33").Value = "  -15.75%  "
$ws.Range("D34").Value = "18.55"
$ws.Range("E34").Value = "  -14.45%  "
$ws.Range("B35").Value = "Monero"
$ws.Range("C35").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D35").Value = "138.12"
$ws.Range("E35").Value = "  -13.83%  "
$ws.Range("B36").Value = "NEARProtocol"
$ws.Range("C36").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D36").Value = "4.13"
$ws.Range("E36").Value = "  -16.25%  "
$ws.Range("D37").Value = "5.36"
$ws.Range("E37").Value = "  -14.81%  "
$ws.Range("D38").Value = "1.20"
$ws.Range("E38").Value = "  -15.59%  "
$ws.Range("D39").Value = "22.85"
$ws.Range("E39").Value = "  -11.32%  "
$ws.Range("D40").Value = "2.897.12"
$ws.Range("E40").Value = "  -9.94%  "
$ws.Range("D41").Value = "0.0612"
$ws.Range("E41").Value = "  -13.02%  "
$ws.Range("D42").Value = "0.999"
$ws.Range("E42").Value = "  -0.20%  "
$ws.Range("D43").Value = "34.97"
$ws.Range("E43").Value = "  -13.98%  "
$ws.Range("E44").Value = "  -14.93%  "
$ws.Range("D45").Value = "0.936"
$ws.Range("E45").Value = "  -13.91%  "
$ws.Range("D46").Value = "1.30"
$ws.Range("E46").Value = "  -11.89%  "
$ws.Range("D47").Value = "3.36"
$ws.Range("E47").Value = "  -15.94%  "
$ws.Range("D48").Value = "2.028.59"
$ws.Range("E48").Value = "  -11.00%  "
$ws.Range("D49").Value = "5.32"
$ws.Range("E49").Value = "  -14.28%  "
$ws.Range("B50").Value = "VeChain"
$ws.Range("C50").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D50").Value = "0.0213"
$ws.Range("E50").Value = "  -10.51%  "
$ws.Range("B51").Value = "InjectiveProtocol"
$ws.Range("C51").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D51").Value = "17.65"
$ws.Range("E51").Value = "  -13.78%  "
